# Update average_county_temperature (column I) and the derived
# worst_ashp_cop / best_ashp_cop (columns N and O) values for the
# rows that were refreshed with NOAA temperature data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (facility_id 1002037, electrified_option 44)
$ws.Range("I4").Value = 17.71296296296294
$ws.Range("N4").Value = 1.872915723725898
$ws.Range("O4").Value = 2.048770944115581

# Row 5 (facility_id 1002037, electrified_option 45)
$ws.Range("I5").Value = 17.71296296296294

# Row 6 (facility_id 1002057, electrified_option 44)
$ws.Range("I6").Value = -1.819444444444444
$ws.Range("N6").Value = 1.572614297115494
$ws.Range("O6").Value = 1.690895540926593

# Row 17 (facility_id 1003352, electrified_option 44)
$ws.Range("I17").Value = -0.763888888888889
$ws.Range("N17").Value = 1.586359976998275
$ws.Range("O17").Value = 1.707009404388715

# Row 20 (facility_id 1003568, electrified_option 44)
$ws.Range("I20").Value = -0.763888888888889
$ws.Range("N20").Value = 1.586359976998275
$ws.Range("O20").Value = 1.707009404388715

# Row 26 (facility_id 1006904, electrified_option 44)
$ws.Range("I26").Value = -1.819444444444444
$ws.Range("N26").Value = 1.572614297115494
$ws.Range("O26").Value = 1.690895540926593

# Row 27 (facility_id 1006904, electrified_option 45)
$ws.Range("I27").Value = -1.819444444444444
